# Apply corrections to the "Nest" sheet (RESULTS_TraitModelTables workbook):
#   - 4 species were missing the nest site trait; this fixes the affected
#     sample sizes / estimates for the "nest site low" and "nest site low
#     (only)" predictor blocks, and removes the now-invalid "nest site high
#     (only)" predictor block entirely (3 models x 3 terms = 9 rows).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Nest")

# --- Row 2-4: predictor_trait "nest site low", model UAI/gls (sample_size 796 -> 792) ---
$ws.Range("E2").Value = 792
$ws.Range("G2").Value = 1.629
$ws.Range("H2").Value = 0.203
$ws.Range("I2").Value = 8.014
$ws.Range("J2").Value = 1.23
$ws.Range("K2").Value = 2.027
$ws.Range("L2").Value = 0.331

$ws.Range("E3").Value = 792
$ws.Range("G3").Value = -0.112
$ws.Range("H3").Value = 0.056
$ws.Range("I3").Value = -2.026
$ws.Range("J3").Value = -0.221
$ws.Range("K3").Value = -0.004
$ws.Range("L3").Value = 0.331

$ws.Range("E4").Value = 792
$ws.Range("H4").Value = 0.023
$ws.Range("I4").Value = -1.207
$ws.Range("J4").Value = -0.071
$ws.Range("L4").Value = 0.331

# --- Row 11-13: predictor_trait "nest site low (only)", model UAI/gls (sample_size 575 -> 571) ---
$ws.Range("E11").Value = 571
$ws.Range("H11").Value = 0.159
$ws.Range("I11").Value = 10.256
$ws.Range("J11").Value = 1.32
$ws.Range("K11").Value = 1.944
$ws.Range("L11").Value = 0.148

$ws.Range("E12").Value = 571
$ws.Range("G12").Value = -0.22
$ws.Range("I12").Value = -3.226
$ws.Range("J12").Value = -0.354
$ws.Range("K12").Value = -0.086
$ws.Range("L12").Value = 0.148

$ws.Range("E13").Value = 571
$ws.Range("G13").Value = -0.024
$ws.Range("H13").Value = 0.022
$ws.Range("I13").Value = -1.09
$ws.Range("J13").Value = -0.066
$ws.Range("L13").Value = 0.148

# --- Remove the "nest site high (only)" predictor block entirely (rows 29-37) ---
# This shifts the remaining "nest safety" and "nest structure (open/enclosed)"
# blocks up by 9 rows (to rows 29-37 and 38-43 respectively), and Excel
# recomputes the shared-strings table (dropping the now-unused
# "nest site high (only)" string) automatically on save.
$ws.Range("A29:M37").EntireRow.Delete()
